# Updates cryptos price (D) and volume(1h) (E) columns per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'26.524.08"
$ws.Range('E2').Value = '  -2.62%  '
$ws.Range('D3').Value = "'1.812.26"
$ws.Range('E3').Value = '  -2.20%  '
$ws.Range('D4').Value = "'1.007"
$ws.Range('E4').Value = '  +0.58%  '
$ws.Range('D5').Value = "'1.005"
$ws.Range('E5').Value = '  +0.45%  '
$ws.Range('D6').Value = "'308.07"
$ws.Range('E6').Value = '  -1.86%  '
$ws.Range('D7').Value = "'0.4554"
$ws.Range('E7').Value = '  -2.12%  '
$ws.Range('D8').Value = "'0.3665"
$ws.Range('E8').Value = '  -1.37%  '
$ws.Range('D9').Value = "'0.07140"
$ws.Range('E9').Value = '  -2.13%  '
$ws.Range('D10').Value = "'0.8796"
$ws.Range('E10').Value = '  -1.12%  '
$ws.Range('D11').Value = "'0.07767"
$ws.Range('E11').Value = '  -1.28%  '
$ws.Range('D12').Value = "'19.37"
$ws.Range('E12').Value = '  -3.73%  '
$ws.Range('D13').Value = "'1.818.59"
$ws.Range('E13').Value = '  -1.83%  '
$ws.Range('E14').Value = '  -2.18%  '
$ws.Range('D15').Value = "'6.371"
$ws.Range('E15').Value = '  -2.29%  '
$ws.Range('D16').Value = "'86.49"
$ws.Range('D17').Value = "'1.007"
$ws.Range('E17').Value = '  +0.54%  '
$ws.Range('D18').Value = "'0.000008599"
$ws.Range('E18').Value = '  -3.61%  '
$ws.Range('D19').Value = "'1.005"
$ws.Range('E19').Value = '  +0.45%  '
$ws.Range('D20').Value = "'26.587.12"
$ws.Range('E20').Value = '  -2.49%  '
$ws.Range('D21').Value = "'14.24"
$ws.Range('E21').Value = '  -3.16%  '
$ws.Range('E22').Value = '  -1.60%  '
$ws.Range('E23').Value = '  -0.64%  '
$ws.Range('D24').Value = "'1.985"
$ws.Range('E24').Value = '  +1.15%  '
$ws.Range('D25').Value = "'151.59"
$ws.Range('E25').Value = '  +0.15%  '
$ws.Range('D26').Value = "'17.94"
$ws.Range('E26').Value = '  -2.41%  '
$ws.Range('D27').Value = "'2.065"
$ws.Range('E27').Value = '  +1.22%  '
$ws.Range('D28').Value = "'112.97"
$ws.Range('E28').Value = '  -2.44%  '
$ws.Range('D29').Value = "'4.856"
$ws.Range('E29').Value = '  -3.59%  '
$ws.Range('D30').Value = "'0.08680"
$ws.Range('E30').Value = '  -1.54%  '
$ws.Range('D31').Value = "'3.059"
$ws.Range('E31').Value = '  -2.58%  '
$ws.Range('D32').Value = "'4.510"
$ws.Range('E32').Value = '  -0.34%  '
$ws.Range('D33').Value = "'0.7345"
$ws.Range('E33').Value = '  -4.26%  '
$ws.Range('D34').Value = "'2.693"
$ws.Range('E34').Value = '  -1.16%  '
$ws.Range('E35').Value = '  -4.10%  '
$ws.Range('D36').Value = "'1.004"
$ws.Range('E36').Value = '  +0.58%  '
$ws.Range('D37').Value = "'1.082"
$ws.Range('E37').Value = '  -2.48%  '
$ws.Range('D38').Value = "'0.01951"
$ws.Range('E38').Value = '  +0.38%  '
$ws.Range('D39').Value = "'0.05115"
$ws.Range('E39').Value = '  -2.08%  '
$ws.Range('D40').Value = "'2.901"
$ws.Range('E40').Value = '  -1.49%  '
$ws.Range('D41').Value = "'6.989"
$ws.Range('E41').Value = '  -0.95%  '
$ws.Range('D42').Value = "'0.5007"
$ws.Range('E42').Value = '  -2.34%  '
$ws.Range('E43').Value = '  -4.16%  '
$ws.Range('E44').Value = '  -3.61%  '
$ws.Range('D45').Value = "'1.005"
$ws.Range('E45').Value = '  +0.47%  '
$ws.Range('D46').Value = "'0.4611"
$ws.Range('E46').Value = '  -3.91%  '
$ws.Range('D47').Value = "'9.963"
$ws.Range('D48').Value = "'100.96"
$ws.Range('E48').Value = '  -1.78%  '
$ws.Range('D49').Value = "'1.593"
$ws.Range('E49').Value = '  -3.18%  '
$ws.Range('E50').Value = '  -3.29%  '
$ws.Range('D51').Value = "'64.46"
$ws.Range('E51').Value = '  -1.50%  '
